$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 42 ("Simulation status and warnings dialog box"),
# pushing that row and everything below it down by one.
$ws.Rows("42:42").Insert()

# Fill in the new task row.
$ws.Range("A42").Value = "Not done"
$ws.Range("B42").Value = "Check SunEye Obstructions and both Solar Pathfinder file imports"
$ws.Range("C42").Value = "Janine"

# Update the active selection to reflect where the user was working.
$ws.Range("A43").Select()
